$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 111; this shifts the existing row 111 (and
# everything below it, through the former row 144) down by one row, which
# matches the diff where rows 111-144 each take on the values that used to
# belong to the row above them, and a brand-new row 145 appears carrying the
# data that used to live in row 144.
$ws.Rows.Item(111).Insert()

# Populate the newly inserted row 111 with the new record.
$ws.Cells.Item(111, 1).Value2 = 5
$ws.Cells.Item(111, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(111, 3).Value2 = "Maule"
$ws.Cells.Item(111, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(111, 4).Value2 = 44543
$ws.Cells.Item(111, 5).Value2 = 7
$ws.Cells.Item(111, 6).Value2 = 100112021
$ws.Cells.Item(111, 7).Value2 = "Ají"
$ws.Cells.Item(111, 8).Value2 = "Americana (o)"
$ws.Cells.Item(111, 9).Value2 = "Primera"
$ws.Cells.Item(111, 10).Value2 = 200
$ws.Cells.Item(111, 11).Value2 = 17000
$ws.Cells.Item(111, 12).Value2 = 17000
$ws.Cells.Item(111, 13).Value2 = 17000
$ws.Cells.Item(111, 14).Value2 = "`$/caja 15 kilos"
$ws.Cells.Item(111, 15).Value2 = "Región del Maule"
$ws.Cells.Item(111, 16).Value2 = 1133
$ws.Cells.Item(111, 17).Value2 = 15
$ws.Cells.Item(111, 18).Value2 = "Hortaliza"
